# TeachingNotes.docx edit script
# 1) Extend the "15 minutes, introductory task" line with extra detail.
# 2) Add the 35 new "ListLabel 323".."ListLabel 357" character styles.

$d = $word.ActiveDocument

# --- 1. Update the "Time Required" line -----------------------------------
$enDash = [char]0x2013
$newText = "15 minutes, introductory task " + $enDash + " for Highers, 1 hour/period - for primary"
$found = $d.Content.Find.Execute("15 minutes, introductory task", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# --- 2. Add the new character styles ---------------------------------------
$newStyles = @(
    @{ Id="ListLabel323"; Name="ListLabel 323"; Ascii="Times New Roman"; Cs=""; Sz=11; U=1 },
    @{ Id="ListLabel324"; Name="ListLabel 324"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel325"; Name="ListLabel 325"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel326"; Name="ListLabel 326"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel327"; Name="ListLabel 327"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel328"; Name="ListLabel 328"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel329"; Name="ListLabel 329"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel330"; Name="ListLabel 330"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel331"; Name="ListLabel 331"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel332"; Name="ListLabel 332"; Ascii="Times New Roman"; Cs=""; Sz=11; U=1 },
    @{ Id="ListLabel333"; Name="ListLabel 333"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel334"; Name="ListLabel 334"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel335"; Name="ListLabel 335"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel336"; Name="ListLabel 336"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel337"; Name="ListLabel 337"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel338"; Name="ListLabel 338"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel339"; Name="ListLabel 339"; Ascii=""; Cs=""; Sz=0; U=1 },
    @{ Id="ListLabel340"; Name="ListLabel 340"; Ascii="Times New Roman"; Cs="OpenSymbol"; Sz=11; U=0 },
    @{ Id="ListLabel341"; Name="ListLabel 341"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel342"; Name="ListLabel 342"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel343"; Name="ListLabel 343"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel344"; Name="ListLabel 344"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel345"; Name="ListLabel 345"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel346"; Name="ListLabel 346"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel347"; Name="ListLabel 347"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel348"; Name="ListLabel 348"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel349"; Name="ListLabel 349"; Ascii="Times New Roman"; Cs="OpenSymbol"; Sz=11; U=0 },
    @{ Id="ListLabel350"; Name="ListLabel 350"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel351"; Name="ListLabel 351"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel352"; Name="ListLabel 352"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel353"; Name="ListLabel 353"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel354"; Name="ListLabel 354"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel355"; Name="ListLabel 355"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel356"; Name="ListLabel 356"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 },
    @{ Id="ListLabel357"; Name="ListLabel 357"; Ascii=""; Cs="OpenSymbol"; Sz=0; U=0 }
)

foreach ($def in $newStyles) {
    $style = $d.Styles.Add($def.Id, 2)
    $style.NameLocal = $def.Name
    $style.QuickStyle = $true
    if ($def.Ascii -ne "") {
        $style.Font.Name = $def.Ascii
    }
    if ($def.Cs -ne "") {
        $style.Font.NameBi = $def.Cs
    }
    if ($def.Sz -gt 0) {
        $style.Font.Size = $def.Sz
    }
    if ($def.U -eq 1) {
        $style.Font.Underline = 0
    }
}

Write-Output "find_result=$found styles_added=$($newStyles.Count)"
